# Fruta / hortaliza, semanal
# Swap the weekly data between rows 2-3 (date 44559) and rows 6-7 (date 44223),
# leaving row 4-5 untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values get swapped between the row pairs
$cols = @("D", "M", "N", "O", "P", "S")

foreach ($col in $cols) {
    # Row 2 <-> Row 6
    $addr2 = $col + "2"
    $addr6 = $col + "6"
    $topVal = $ws.Range($addr2).Value2
    $botVal = $ws.Range($addr6).Value2
    $ws.Range($addr2).Value2 = $botVal
    $ws.Range($addr6).Value2 = $topVal

    # Row 3 <-> Row 7
    $addr3 = $col + "3"
    $addr7 = $col + "7"
    $topVal2 = $ws.Range($addr3).Value2
    $botVal2 = $ws.Range($addr7).Value2
    $ws.Range($addr3).Value2 = $botVal2
    $ws.Range($addr7).Value2 = $topVal2
}
